$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Worksheet data changes -------------------------------------------------
# "actual remaining story points" (column E) now has data for the first
# 9 days of the sprint (all still at 50 story points remaining).
$ws.Range("E7:E15").Value = 50

# "ideal remaining story points" (column F) burndown formula now starts
# from 50 points instead of 18. Re-entering the formula across F7:F24
# recreates it as a shared formula, same as the original authoring.
$ws.Range("F7:F24").Formula = "=50 - (50*(D7-1)/17)"

# Selection moved from F24 to F3, and the view scrolled back to show
# column C (the start of the table) - selecting F3 mirrors that.
$ws.Range("F3").Select()

# --- Chart changes -----------------------------------------------------
$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart

# Value axis maximum raised from 18 to 50 to match the new scale.
$valAx = $chart.Axes(2)
$valAx.MaximumScale = 50

# The chart's anchor (twoCellAnchor) moved/resized on the sheet: it now
# starts near column G / row 2 and ends near column T / row 39, instead
# of the previous column G / row 8 -> column R / row 32 placement.
# Translate the target cell+offset anchor into absolute points (the units
# ChartObject.Left/Top/Width/Height use) via the corresponding cell
# geometry already on the sheet.
$fromLeft = $ws.Range("G1").Left + (169067 / 12700)
$fromTop  = $ws.Range("A2").Top + (171450 / 12700)
$toLeft   = $ws.Range("T1").Left + (390525 / 12700)
$toTop    = $ws.Range("A39").Top + (123825 / 12700)

$chartObj.Left = $fromLeft
$chartObj.Top = $fromTop
$chartObj.Width = $toLeft - $fromLeft
$chartObj.Height = $toTop - $fromTop
